# Delete column A entirely; remaining columns (B:F) shift left to become A:E,
# keeping each cell's own formatting/content intact.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A:A").Delete()
